# Delete column E ("reviews_count") entirely, shifting columns F:K left to E:J.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete()
